# Daily IST report: add CSV/MD/XLSX
# Adds a new daily-submission column for 2026-02-28 into the
# "daily_counts" sheet (inserted before the existing total_files /
# unique_days summary columns), then recomputes the two summary
# columns for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column at N (pushes total_files/unique_days from N/O to O/P)
# ---------------------------------------------------------------------------
$refWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N:N").Insert()
$ws.Columns("N").ColumnWidth = $refWidth

# ---------------------------------------------------------------------------
# 2. New header cell N1 = "2026-02-28" (plain text, same style as the other
#    date headers - copy format only from M1 so the bold/centered xf is
#    reused instead of minting a brand-new style).
# ---------------------------------------------------------------------------
$ws.Range("N1").NumberFormat = "@"
$ws.Range("N1").Value = "2026-02-28"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill in the new 2026-02-28 submission flag for every data row. Rows not
#    listed below submitted nothing that day and get an explicit literal 0
#    (matching the source data's style of writing 0s rather than leaving the
#    cell blank).
# ---------------------------------------------------------------------------
$rowsWithSubmission = @(4,6,9,10,12,18,19,31,34,42,45,46,49,53,60,70,76,81,103)
for ($r = 2; $r -le 109; $r++) {
    $ws.Cells.Item($r, 14).Value = 0
}
foreach ($r in $rowsWithSubmission) {
    $ws.Cells.Item($r, 14).Value = 1
}

# ---------------------------------------------------------------------------
# 4. Recompute total_files (O, col 15) = SUM(D:N) and
#               unique_days  (P, col 16) = COUNT of D:N cells > 0
#    for every data row, now that the daily matrix spans D..N.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 109; $r++) {
    $total = 0
    $days = 0
    for ($c = 4; $c -le 14; $c++) {
        $val = $ws.Cells.Item($r, $c).Value2
        $total = $total + $val
        if ($val -gt 0) {
            $days = $days + 1
        }
    }
    $ws.Cells.Item($r, 15).Value = $total
    $ws.Cells.Item($r, 16).Value = $days
}
